# Rúbrica edit: "commit rubrica 4a entrega2"
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append ". " right after "...capa de sistemas" (end of that
#    paragraph, before the paragraph mark).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "documentos del negocio, reportes etc.. que tendrán relaciones con set de datos en la capa de sistemas",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(". ")

# ------------------------------------------------------------------
# 2) Append a new run after "... para que puedan realizar análisis de
#    impacto." and then add four brand new paragraphs describing the
#    minimum component counts per layer.
#
#    Each of the new paragraphs is built from several small runs (to
#    match how the source commit recorded them). Plain sequential
#    Range.InsertAfter calls get silently coalesced into one run when
#    they share formatting, so the runs are typed with TrackRevisions
#    on (which keeps every insertion as its own run/w:ins) and then
#    accepted in one shot — the accepted result keeps the run
#    boundaries while dropping all the revision markup.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "para que puedan realizar análisis de impacto.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)

$d.TrackRevisions = $true

$rng.InsertAfter(" Luego debe realizar análisis de impacto. ")
$rng.Collapse(0)

# --- new paragraph: capa de motivación ---
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Cantidad de componentes")
$rng.Collapse(0)
$rng.InsertAfter(" mínimo")
$rng.Collapse(0)
$rng.InsertAfter(" para la capa de motivación: ")
$rng.Collapse(0)
$rng.InsertAfter("1")
$rng.Collapse(0)
$rng.InsertAfter("2")
$rng.Collapse(0)
$rng.InsertAfter(" de diferente tipo")
$rng.Collapse(0)
$rng.InsertAfter(": assesment, ")
$rng.Collapse(0)
$rng.InsertAfter("drivers, objetivos, metas, requerimientos, stakeholders")
$rng.Collapse(0)

# --- new paragraph: capa de negocio ---
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Cantidad de componentes mínimo para la capa de ")
$rng.Collapse(0)
$rng.InsertAfter("negocio")
$rng.Collapse(0)
$rng.InsertAfter(": 1")
$rng.Collapse(0)
$rng.InsertAfter("5")
$rng.Collapse(0)
$rng.InsertAfter(" de diferente tipo")
$rng.Collapse(0)
$rng.InsertAfter(", deben haber roles, actores, procesos, funciones, catálogo")
$rng.Collapse(0)
$rng.InsertAfter("s y documentos")
$rng.Collapse(0)

# --- new paragraph: capa de sistemas ---
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Cantidad de compone")
$rng.Collapse(0)
$rng.InsertAfter("ntes de la capa de sistemas: los necesarios para representar la arquitectura que implementó")
$rng.Collapse(0)

# --- new paragraph: capa de infraestructura ---
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Cantidad de componentes mínimo de la capa de infraestructura: 1")
$rng.Collapse(0)
$rng.InsertAfter("2")
$rng.Collapse(0)
$rng.InsertAfter(". Deben ser coherentes para que allí esté")
$rng.Collapse(0)
$rng.InsertAfter(" implementada toda la solución de la capa superior.")
$rng.Collapse(0)

$d.TrackRevisions = $false
$d.AcceptAllRevisions() | Out-Null

# ------------------------------------------------------------------
# 3) Move the "lastRenderedPageBreak" rendering hint down one table
#    row: it now falls on the "Tiene elementos que describen la
#    estrategia" run instead of the "Hay roles y actores..." run,
#    because the new content inserted above shifted the page break.
#    There's no Word-OM property for this low-level layout marker, so
#    it's done with a scoped Range.InsertXML that rewrites just the
#    owning paragraph, preserving its pPr/run text.
# ------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t -like "*Tiene elementos que describen la estrategia*") {
        $full = $para.Range
        $target = $d.Range($full.Start, $full.End - 1)
        $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:ind w:left="708"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Tiene elementos que describen la estrategia</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $target.InsertXML($xmlFrag) | Out-Null
    }
    elseif ($t -like "*Hay roles y actores de la organización para este componente*") {
        $full = $para.Range
        $target = $d.Range($full.Start, $full.End - 1)
        $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:ind w:left="708"/></w:pPr><w:r><w:t>Hay roles y actores de la organización para este componente</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $target.InsertXML($xmlFrag) | Out-Null
    }
}
